$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlLeft constant
$xlLeft = -4131

# Fill in "x" marks for H8, H9, H10, H13 (column H was left blank before)
$ws.Range("H8").Value = "x"
$ws.Range("H9").Value = "x"
$ws.Range("H10").Value = "x"
$ws.Range("H13").Value = "x"

# Row 20: add "Client has" note in H20 and start the code snippet list in I20
$ws.Range("H20").Value = "Client has"

$ws.Range("I20").Value = "self.addMethodToType(clMethod)"
$ws.Range("I20").WrapText = $false
$ws.Range("I20").HorizontalAlignment = $xlLeft

$ws.Range("I21").Value = "self.addMethodToTyoeFromDB(iMethodId)"
$ws.Range("I21").HorizontalAlignment = $xlLeft

$ws.Range("I22").Value = "self.addPropertyToType(clProperty)"
$ws.Range("I22").HorizontalAlignment = $xlLeft

$ws.Range("I23").Value = "self.addEventToType(clEvent)"
$ws.Range("I23").HorizontalAlignment = $xlLeft

$ws.Range("I24").HorizontalAlignment = $xlLeft

$ws.Range("I25").Value = "These are called from NewMethod, SearchForMethod, NewProperty and NewEvent dialogs."
$ws.Range("I25").HorizontalAlignment = $xlLeft

# Move the active selection to H14
$ws.Range("H14").Select()
